# Weekly update of Fruit/Vegetable price data.
# Updates columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion), O (Origen),
# P (Precio $/Kg) and Q (Kg o Unidades) for rows 3-9,11-13 as the weekly
# dataset rotates through the records (row 10 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $D, $J, $K, $L, $M, $N, $O, $P, $Q) {
    $ws.Range("D$Row").Value = $D
    $ws.Range("J$Row").Value = $J
    $ws.Range("K$Row").Value = $K
    $ws.Range("L$Row").Value = $L
    $ws.Range("M$Row").Value = $M
    $ws.Range("N$Row").Value = $N
    $ws.Range("O$Row").Value = $O
    $ws.Range("P$Row").Value = $P
    $ws.Range("Q$Row").Value = $Q
}

Set-Row 3  44315 25 10000 10000 10000 "$/caja 60 unidades" "Provincia de Limarí"   167 60
Set-Row 4  44312 30 10000 10000 10000 "$/caja 60 unidades" "Provincia de Limarí"   167 60
Set-Row 5  44284 35 10000 10000 10000 "$/caja 60 unidades" "Provincia de Limarí"   167 60
Set-Row 6  44405 45 9000  9000  9000  "$/caja 50 unidades" "Provincia de Quillota" 180 50
Set-Row 7  44277 25 10000 10000 10000 "$/caja 60 unidades" "Provincia de Limarí"   167 60
Set-Row 8  44186 15 7000  7000  7000  "$/caja 60 unidades" "Provincia de Limarí"   117 60
Set-Row 9  44243 80 10000 11000 10375 "$/caja 60 unidades" "Provincia de Quillota" 173 60
Set-Row 11 44585 30 11000 11000 11000 "$/caja 60 unidades" "Provincia de Limarí"   183 60
Set-Row 12 44179 15 7000  7000  7000  "$/caja 60 unidades" "Provincia de Limarí"   117 60
Set-Row 13 44291 20 9000  9000  9000  "$/caja 60 unidades" "Provincia de Limarí"   150 60
